$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.248.20'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '1.862.69'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.94'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9992'
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4700'
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('E8').Value = '  +2.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06573'
$ws.Range('E9').Value = '  +0.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.87'
$ws.Range('E10').Value = '  +1.93%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08005'
$ws.Range('E11').Value = '  +1.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '97.76'
$ws.Range('E12').Value = '  -0.53%  '
$ws.Range('D13').Value = '1.858.15'
$ws.Range('E13').Value = '  -0.68%  '
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6796'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '271.52'
$ws.Range('E16').Value = '  -2.48%  '
$ws.Range('D17').Value = '30.226.67'
$ws.Range('E18').Value = '  +7.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007653'
$ws.Range('E19').Value = '  +4.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9994'
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('D21').Value = '2.104.47'
$ws.Range('E21').Value = '  -0.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9992'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.238'
$ws.Range('E23').Value = '  -4.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.186'
$ws.Range('E24').Value = '  +0.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '167.13'
$ws.Range('E25').Value = '  +1.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.202'
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.97'
$ws.Range('E27').Value = '  -1.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.952'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.371'
$ws.Range('E29').Value = '  -1.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09949'
$ws.Range('E30').Value = '  +2.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.341'
$ws.Range('E31').Value = '  -1.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.466'
$ws.Range('E32').Value = '  -0.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.046'
$ws.Range('E33').Value = '  -1.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04710'
$ws.Range('E34').Value = '  -0.21%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.125'
$ws.Range('E35').Value = '  -0.70%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7044'
$ws.Range('E36').Value = '  -0.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.716'
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01881'
$ws.Range('E38').Value = '  +0.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.602'
$ws.Range('E39').Value = '  +2.48%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.346'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '73.50'
$ws.Range('E41').Value = '  -1.89%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.943'
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '103.95'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8402'
$ws.Range('E44').Value = '  -1.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9984'
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4152'
$ws.Range('E46').Value = '  -1.06%  '
$ws.Range('E47').Value = '  -2.38%  '
$ws.Range('E48').Value = '  -0.77%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '931.89'
$ws.Range('E49').Value = '  -1.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.20'
$ws.Range('E50').Value = '  -0.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05666'
$ws.Range('E51').Value = '  +0.43%  '
